$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Matthew Wolz"
$ws.Range("B6").Value = "Greg s"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "3/24/2025"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").Value = "MW"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = "Daily Guest Pass"
$ws.Range("H6").Value = 3
